$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths -------------------------------------------------------
# Column A -> 13.0 chars, Column C -> 32.14 chars. Column B and D:Y stay at
# the sheet's existing default (8.71) and are left untouched.
$ws.Columns.Item(1).ColumnWidth = 12.17
$ws.Columns.Item(3).ColumnWidth = 31.307

# --- Header row (row 1) ---------------------------------------------------
$ws.Range("A1").Value = "username"
$ws.Range("B1").Value = "gender"
$ws.Range("C1").Value = "gmail"
$ws.Range("D1").Value = "status"

# --- Row 2 -----------------------------------------------------------------
$ws.Range("A2").Value = "Thành Côngdncv"
$ws.Range("B2").Value = "f"
$ws.Range("C2").Value = "nguyenthanhcong.dn.cv@gmail.com"
$ws.Range("D2").ClearContents()

# --- Row 3 -----------------------------------------------------------------
$ws.Range("A3").Value = "Thành Công120620"
$ws.Range("B3").Value = "m"
$ws.Range("C3").Value = "nguyenthanhcong120620@gmail.com"
$ws.Range("D3").ClearContents()

# --- Row 4 -----------------------------------------------------------------
$ws.Range("A4").Value = "Thành Côngivcgroup"
$ws.Range("B4").Value = "fm"
$ws.Range("C4").Value = "cong.nt.ivcgroup@gmail.com"
$ws.Range("D4").ClearContents()

# --- Row 5 -----------------------------------------------------------------
$ws.Range("A5").Value = "Thành Côngdev126"
$ws.Range("B5").Value = "m"
$ws.Range("C5").Value = "nguyenthanhcong.dev.126@gmail.com"
$ws.Range("D5").ClearContents()
